$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing entry for jagadeesh k.docx) -- new recommended job: Basil Technologies
$ws.Range("E2").Value = "https://www.mycareersfuture.gov.sg/job/information-technology/senior-developer-basil-technologies-0ca878162e09111d2be8ee4a64644aec"
$ws.Range("F2").Value = "MCF-2022-0296303"
$ws.Range("G2").Value = "BASIL TECHNOLOGIES PTE. LTD."
$ws.Range("H2").Value = "Senior Developer "
$ws.Range("I2").Value = "56% skills matched"
$ws.Range("J2").Value = "SOAP, MySQL, Configuration Management, API, Databases, J2EE, Hibernate, "
$ws.Range("K2").Value = "Contract"
$ws.Range("L2").Value = "Senior Executive"
$ws.Range("M2").Value = "Information Technology"
$ws.Range("N2").Value = "Islandwide"
$ws.Range("O2").Value = "6 years exp"
$ws.Range("P2").Value = '$6,000 to $8,000 Monthly'

# Row 3 (existing entry for Jennifer M. Conte.docx) -- new recommended job: Marquee Semiconductor
$ws.Range("E3").Value = "https://www.mycareersfuture.gov.sg/job/engineering/senior-software-engineers-marquee-semiconductor-singapore-c62f0f74f24a099da4b542e215d701ac"
$ws.Range("F3").Value = "MCF-2022-0282201"
$ws.Range("G3").Value = "MARQUEE SEMICONDUCTOR  SINGAPORE PTE. LTD."
$ws.Range("H3").Value = "Senior Software Engineers"
$ws.Range("I3").Value = "30% skills matched"
$ws.Range("J3").Value = "Version Control, Autonomy, AngularJS, MySQL, Scripting, Information Technology, OpenCL, Requirements Analysis, Python, Computer Architecture, ClearCase, Debugging, Databases, Software Development, "
$ws.Range("K3").Value = "Full Time"
$ws.Range("L3").Value = "Middle Management"
$ws.Range("M3").Value = "Engineering"
$ws.Range("N3").Value = "Islandwide"
$ws.Range("O3").Value = "1 year exp"
$ws.Range("P3").Value = '$4,000 to $6,000 Monthly'

# Row 4 (existing entry for SUNITHA Project Manager (1).docx) -- new recommended job: The Supreme HR Advisory
$ws.Range("E4").Value = "https://www.mycareersfuture.gov.sg/job/information-technology/business-analyst-5-days-java-2894-supreme-hr-advisory-81fc54813d31f893d647d8ddfa32cb8f"
$ws.Range("F4").Value = "MCF-2022-0332440"
$ws.Range("G4").Value = "THE SUPREME HR ADVISORY PTE. LTD."
$ws.Range("H4").Value = "Business Analyst [5 days| Java] 2894"
$ws.Range("I4").Value = "50% skills matched"
$ws.Range("J4").Value = "Requirements Gathering, Microsoft Excel, Business Analysis, Communication Skills, Banking, Web Applications, Team Player, Business Requirements, "
$ws.Range("K4").Value = "Permanent"
$ws.Range("L4").Value = "Junior Executive"
$ws.Range("M4").Value = "Information Technology"
$ws.Range("N4").Value = "Islandwide"
$ws.Range("O4").Value = "1 year exp"
$ws.Range("P4").Value = '$3,200 to $5,000 Monthly'
